# TABELA REFERENCIAL GTI 01 - add three new line items to the "gti" sheet
# (TAMPA CEGA 4X2 / CABO FIBRA ÓPTICA / CORDÃO ÓPTICO MULTIMODO LC) just
# above the closing "*OBS" / total row, pushing that row from 88 down to 91.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows right before the closing row (row 88), which
# shifts the existing closing row down to row 91.
$ws.Rows("88:90").Insert()

# Clone the formatting of the row directly above (row 87, the last data
# row) onto the three new rows so borders/fonts/number formats match the
# rest of the table exactly.
$ws.Range("A87:F87").Copy()
$ws.Range("A88:F90").PasteSpecial(-4122)

# New items to add, in order.
$items = @("TAMPA CEGA 4X2", "CABO FIBRA ÓPTICA", "CORDÃO ÓPTICO MULTIMODO LC")

for ($i = 0; $i -lt $items.Length; $i++) {
    $r = 88 + $i
    $ws.Cells.Item($r, 1).Value = "PNCP"
    $ws.Cells.Item($r, 2).Value = "PNCP"
    $ws.Cells.Item($r, 3).Value = $items[$i]
    $ws.Cells.Item($r, 4).Value = "NUMERO DE CONTROLE 01 "
    $ws.Cells.Item($r, 5).Value = "NUMERO DE CONTROLE 02 "
    $ws.Cells.Item($r, 6).Value = "NUMERO DE CONTROLE 03"
}

# Match the saved view state: selection on D89:F90 (scrolled near the top
# of the sheet, around row 9).
$ws.Activate() | Out-Null
$ws.Range("D89:F90").Select() | Out-Null
